$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph: "Ultimately, our problem is a classification problem. ... sliding
# or rotating state. " -> append the new trailing sentence about the 10-second
# prediction window, and carry the (hidden) _GoBack bookmark along to the new
# end of the paragraph, the way Word does when new text is typed at the spot
# it's tracking.
# ---------------------------------------------------------------------------
$targetText = "Ultimately, our problem is a classification problem."
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith($targetText)) {
        $para = $cand
        break
    }
}
if ($para -eq $null) {
    throw "Could not locate the 'Ultimately, our problem...' paragraph"
}

$paraIndex = $para.Index
$insertPos = $para.Range.End - 1
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter("At this point, we are only predicting the very next data point which is 10 seconds out.")

# Re-resolve the paragraph (its end moved) and find the new end-of-text spot,
# immediately before the paragraph mark.
$para = $d.Paragraphs.Item($paraIndex)
$newEndPos = $para.Range.End - 1

# Relocate the _GoBack bookmark there. A collapsed range landing exactly on a
# paragraph-mark position trips up Bookmarks.Add in this host, so nudge a
# placeholder character in first, anchor the bookmark ahead of it, then drop
# the placeholder back out.
if ($d.Bookmarks.Exists("_GoBack")) {
    $placeholder = $d.Range($newEndPos, $newEndPos)
    $placeholder.InsertAfter("Z")
    $bmRange = $d.Range($newEndPos, $newEndPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    $cleanup = $d.Range($newEndPos, $newEndPos + 1)
    $cleanup.Text = ""
}

# ---------------------------------------------------------------------------
# Paragraph: "All 3 of our models performed well ... almost 99% accuracy." ->
# append the follow-up sentence comparing logistic regression's speed.
# ---------------------------------------------------------------------------
$targetText2 = "All 3 of our models performed well"
$para2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith($targetText2)) {
        $para2 = $cand
        break
    }
}
if ($para2 -eq $null) {
    throw "Could not locate the 'All 3 of our models...' paragraph"
}

$insertPos2 = $para2.Range.End - 1
$ins2 = $d.Range($insertPos2, $insertPos2)
$ins2.InsertAfter(" Logistic regression performed almost as well as SVM and ran much quicker.")

Write-Output "done"
